$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the province names in row 53 and 54 (Huelva <-> Huesca)
$ws.Range("A53").Value = "Huesca"
$ws.Range("A54").Value = "Huelva"

# Swap the "Casos activos" values between the two rows
$ws.Range("C53").Value = 0
$ws.Range("C54").Value = 72

# Update the death counts (column E) for several rows
$ws.Range("E32").Value = 8
$ws.Range("E56").Value = 8
$ws.Range("E57").Value = 8
$ws.Range("E65").Value = 8

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 10:46"
